# The original sheet had two rows in column A:
#   A1 = 0                (bold text, centered, thin border)
#   A2 = the Python-dict-style "questions = [...]" string (default style)
#
# The edit removes the standalone "0" row and reformats the remaining
# questions string as pretty-printed JSON-like text (also fixing the
# trailing space in the second question's "score" value). After the
# edit there is a single row containing only the reformatted string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop row 1 (the "0" cell) - this shifts the questions string up to A1
# and drops the bold/bordered formatting that only applied to row 1.
$ws.Rows(1).Delete()

# Replace the cell's text with the reformatted (pretty-printed) version.
$ws.Range("A1").Value = 'questions = [
    {
        "title": "You are developing an advanced deep-learning model for a security system that requires real-time facial recognition. The model needs to maintain high accuracy in various lighting conditions, which is challenging due to the complexity of facial features and the variability in image quality.Which of the following architectural modifications should you prioritize?",
        "ques_type": 2,
        "options": [
            "Implement a deep convolutional neural network (CNN) with residual connections.",
            "Use a shallow convolutional neural network (CNN) with minimal layers.",
            "Apply a recurrent neural network (RNN) layer.",
            "Opt for a generative adversarial network (GAN)."
        ],
        "score": "Implement a deep convolutional neural network (CNN) with residual connections."
    },
    {
        "title": "You''re training a sophisticated neural network for natural language processing tasks. The network''s convergence is slower than expected, potentially delaying the project''s timeline and affecting performance.Which optimization algorithm should you use?",
        "ques_type": 2,
        "options": [
            "Adam optimizer ",
            "Stochastic gradient descent ",
            "RMSprop",
            "Adagrad optimizer"
        ],
        "score": "Adam optimizer"
    },
    {
        "title": "Your team needs to enhance the computer vision system of an autonomous vehicle. The focus is on developing a robust neural network model that can accurately detect and classify road signs, especially in varying weather conditions.Which of the following neural network approaches should you use?",
        "ques_type": 2,
        "options": [
            "Convolutional neural network (CNN) with data augmentation",
            "Convolutional neural network (CNN) without data augmentation",
            "Recurrent neural network (RNN) without data augmentation",
            "Recurrent neural network (RNN) with data augmentation"
        ],
        "score": "Convolutional neural network (CNN) with data augmentation"
    },
    {
        "title": "To improve your company''s chatbot, you''re exploring the latest advancements in neural network models. The goal is to enhance the chatbot''s ability to provide contextually relevant responses, a key factor in user satisfaction and engagement.Which neural network models should you use?",
        "ques_type": 2,
        "options": [
            "Transformer-based models",
            "Traditional long short-term memory (LSTM) networks ",
            "Convolutional neural networks (CNNs)",
            "Autoencoders"
        ],
        "score": "Transformer-based models"
    }
]'

# Re-fit the row height to the default so no stale explicit row height
# (left over from the previous long single-line text) is persisted.
$ws.Rows(1).AutoFit()
